$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits: shuffle the raw Raman-shift data points ---------------
# The first two raw values (old B2, old B3) are demoted to plain numbers
# and moved down into B4/B5 (no longer part of the live average/stdev).
# The last two raw values (old B4, old B5) become the new live formulas
# in B2/B3 that feed the summary statistics.

$ws.Range("B2").Formula = "=-0.05-0.11"
$ws.Range("B3").Formula = "=0.01-0.14"

$ws.Range("B4").Value = -0.13
$ws.Range("B5").Value = -0.19

# Formatting applied to the now-static data points B4/B5
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B4").VerticalAlignment = -4108  # xlCenter (vertical)

$ws.Range("B5").Font.Name = "Calibri"
$ws.Range("B5").NumberFormat = "0.00"

# Summary formulas now reference the new B4:B5 (and the quirky B5:B6 range
# for the std-dev, exactly as re-authored)
$ws.Range("B6").Formula = "=AVERAGE(B4:B5)"
$ws.Range("B7").Formula = "=STDEV.S(B5:B6)"

# The D column helper cells (blank, style-only) are no longer used
$ws.Range("D6").Clear()
$ws.Range("D7").Clear()

# --- View / selection state ---------------------------------------------
$ws.Range("F11").Select()

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait
